# Update "想去人数" (want-to-go count) values in the F column across sheets.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 96
$ws1.Range("F3").Value = 1225
$ws1.Range("F4").Value = 864
$ws1.Range("F5").Value = 889
$ws1.Range("F6").Value = 1625
$ws1.Range("F7").Value = 342
$ws1.Range("F8").Value = 1098
$ws1.Range("F11").Value = 228
$ws1.Range("F12").Value = 71
$ws1.Range("F13").Value = 573
$ws1.Range("F14").Value = 100
$ws1.Range("F19").Value = 33
$ws1.Range("F20").Value = 611
$ws1.Range("F21").Value = 602
$ws1.Range("F22").Value = 90
$ws1.Range("F24").Value = 814
$ws1.Range("F25").Value = 274
$ws1.Range("F27").Value = 218
$ws1.Range("F29").Value = 387

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 77
$ws2.Range("F8").Value = 102

# Sheet "本地生活" (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 284

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 284
$ws4.Range("F3").Value = 96
$ws4.Range("F4").Value = 1225
$ws4.Range("F5").Value = 864
$ws4.Range("F6").Value = 889
$ws4.Range("F7").Value = 1625
$ws4.Range("F8").Value = 342
$ws4.Range("F9").Value = 1098
$ws4.Range("F12").Value = 228
$ws4.Range("F13").Value = 71
$ws4.Range("F14").Value = 573
$ws4.Range("F15").Value = 100
$ws4.Range("F25").Value = 33
$ws4.Range("F26").Value = 611
$ws4.Range("F27").Value = 602
$ws4.Range("F28").Value = 90
$ws4.Range("F30").Value = 814
$ws4.Range("F31").Value = 274
$ws4.Range("F32").Value = 77
$ws4.Range("F34").Value = 218
$ws4.Range("F36").Value = 102
$ws4.Range("F37").Value = 102
$ws4.Range("F40").Value = 387

$wb.Save()
